# PlayerPerformance_4714.xlsx: add "Player Info" sheet (front) and
# "ODI Batting Extra" sheet (end); rename MATCH_CARD_LINK -> MATCH_CODE on
# the batting/bowling sheets and collapse the link URL down to the bare
# numeric match code; drop the stray empty INNING_NUMBER cells on ODI Batting.
#
# NOTE: worksheet COM handles returned by Worksheets.Item(...) resolve by
# *position*, not stable identity, in this host - so every lookup re-fetches
# the sheet by name right before it is used, rather than caching a reference
# across a structural change (sheet insertion) that could shift indices.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122
$xlDown = -4121

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $playerInfoHeaders.Length; $col++) {
    $playerInfo.Cells.Item(1, $col).Value = $playerInfoHeaders[$col - 1]
}

# match header formatting used on the other sheets (bold, centered, bordered)
$wb.Worksheets.Item("ODI Batting").Cells.Item(1, 1).Copy()
$playerInfo.Range($playerInfo.Cells.Item(1, 1), $playerInfo.Cells.Item(1, 4)).PasteSpecial($xlPasteFormats)

$playerInfoRow = @("4714", "Mujeeb Ur Rahman", "Right Handed", "Right Arm Off Break")
for ($col = 1; $col -le $playerInfoRow.Length; $col++) {
    $cell = $playerInfo.Cells.Item(2, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $playerInfoRow[$col - 1]
}

# ---------------------------------------------------------------------------
# 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code,
#    drop empty INNING_NUMBER cells
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $battingSheet.Cells.Item(1, 1).End($xlDown).Row
for ($row = 2; $row -le $battingLastRow; $row++) {
    $linkCell = $battingSheet.Cells.Item($row, 4)
    $link = $linkCell.Value2
    if ($link) {
        $code = $link.ToString().Split("=")[-1]
        $linkCell.NumberFormat = "@"
        $linkCell.Value2 = $code
    }

    $inningCell = $battingSheet.Cells.Item($row, 2)
    if ([string]::IsNullOrEmpty($inningCell.Value2)) {
        $inningCell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.Cells.Item(1, 1).End($xlDown).Row
for ($row = 2; $row -le $bowlingLastRow; $row++) {
    $linkCell = $bowlingSheet.Cells.Item($row, 2)
    $link = $linkCell.Value2
    if ($link) {
        $code = $link.ToString().Split("=")[-1]
        $linkCell.NumberFormat = "@"
        $linkCell.Value2 = $code
    }
}

# ---------------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet, appended after "ODI Bowling"
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Bowling"))
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $extraHeaders[$col - 1]
}
$wb.Worksheets.Item("ODI Batting").Cells.Item(1, 1).Copy()
$extra.Range($extra.Cells.Item(1, 1), $extra.Cells.Item(1, 6)).PasteSpecial($xlPasteFormats)

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4340", 11, "1", "0", "3.08%", "NO"),
    @("4348", 11, "0", "0", "", "NO"),
    @("4377", 10, "0", "0", "", "NO"),
    @("4378", 10, "", "", "", "NO"),
    @("4379", 10, "0", "0", "0.35%", "NO"),
    @("4444", 10, "", "", "", "NO"),
    @("4446", 10, "2", "0", "6.77%", "NO"),
    @("4448", "", "", "", "", "NO"),
    @("4525", 9, "", "", "", "NO"),
    @("4528", 8, "", "", "", "NO"),
    @("4530", 9, "0", "0", "", "NO"),
    @("4537", "", "", "", "", "NO"),
    @("4538", 9, "", "", "", "NO"),
    @("4539", 9, "", "", "", "NO"),
    @("4582", 9, "", "", "", "NO"),
    @("4585", 9, "", "", "", "NO"),
    @("4588", "", "", "", "", "NO"),
    @("4671", 9, "0", "0", "", "NO"),
    @("4674", "", "", "", "", "NO"),
    @("4675", "", "", "", "", "")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = $i + 2
    $values = $extraRows[$i]

    $codeCell = $extra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $values[0]

    $posCell = $extra.Cells.Item($r, 2)
    if ($values[1] -ne "") {
        $posCell.Value = $values[1]
    }

    for ($col = 3; $col -le 5; $col++) {
        $v = $values[$col - 1]
        if ($v -ne "") {
            $c = $extra.Cells.Item($r, $col)
            $c.NumberFormat = "@"
            $c.Value = $v
        }
    }

    $momCell = $extra.Cells.Item($r, 6)
    if ($values[5] -ne "") {
        $momCell.Value = $values[5]
    }
}
